$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.6945316142291066
$ws.Range("C2").Value = 0.2103628180572122
$ws.Range("D2").Value = -0.5445742289801792
$ws.Range("E2").Value = 1.649195601823524
$ws.Range("F2").Value = 1.939470983457267
$ws.Range("B3").Value = 0.05114009865011426
$ws.Range("C3").Value = 1.690456444252163
$ws.Range("D3").Value = 2.291947805260176
$ws.Range("E3").Value = 2.322941629478642
$ws.Range("F3").Value = -1.830493192148062
$ws.Range("G3").Value = 0.4490046210795488
$ws.Range("H3").Value = -1.754915860385852
$ws.Range("B4").Value = 1.744455382806309
$ws.Range("C4").Value = 1.561954287081898
$ws.Range("D4").Value = -2.344749236599752
$ws.Range("E4").Value = 0.1769446262318977
$ws.Range("F4").Value = -1.956889776802553
$ws.Range("B5").Value = -2.144234228514782
$ws.Range("C5").Value = 0.5280052638876265
$ws.Range("D5").Value = -1.593411734382765
$ws.Range("E5").Value = 1.049026178229147
$ws.Range("F5").Value = -0.7824510563527671
$ws.Range("G5").Value = -0.1441541880888198
$ws.Range("H5").Value = -0.3865776622619923
$ws.Range("B6").Value = -1.301036791611902
$ws.Range("C6").Value = 1.28984165056396
$ws.Range("D6").Value = -0.6431513272560402
$ws.Range("E6").Value = -0.1433097936191611
$ws.Range("F6").Value = -0.3942967594327705
$ws.Range("B7").Value = 0.264747157129269
$ws.Range("C7").Value = 0.2184107288059266
$ws.Range("D7").Value = -0.4791188345420213
$ws.Range("E7").Value = -0.3240821062975808
$ws.Range("F7").Value = 0.06521323469150679
$ws.Range("G7").Value = 0.01983498341556333
$ws.Range("H7").Value = 0.5402734837222156
$ws.Range("B8").Value = -0.5072360807453689
$ws.Range("C8").Value = -0.4720094128878891
$ws.Range("D8").Value = [double]"-9.71445146547012E-17"
$ws.Range("E8").Value = -0.0794622193825073
$ws.Range("F8").Value = 0.5
$ws.Range("B9").Value = 0.1999999999999999
$ws.Range("C9").Value = -0.0477857303545277
$ws.Range("D9").Value = 0.4639869465896599
$ws.Range("E9").Value = -0.8
$ws.Range("F9").Value = -0.6749228041047655
$ws.Range("G9").Value = -0.2410120038038291
$ws.Range("H9").Value = 0.7358685032995328
$ws.Range("B10").Value = 0.1
$ws.Range("C10").Value = -1.2
$ws.Range("D10").Value = -0.8999999999999999
$ws.Range("E10").Value = -0.2872048886675244
$ws.Range("F10").Value = 0.7027743565091521
$ws.Range("B11").Value = -1.1
$ws.Range("C11").Value = -0.5
$ws.Range("D11").Value = 0.6000000000000001
$ws.Range("F11").Value = -0.3999999999999999
$ws.Range("G11").Value = 0.4079191665375851
$ws.Range("H11").Value = -0.7011422248794972
$ws.Range("B12").Value = 0.3
$ws.Range("C12").Value = [double]"2.775557561562891E-17"
$ws.Range("D12").Value = -0.2999999999999999
$ws.Range("E12").Value = 0.3912510204569938
$ws.Range("F12").Value = -0.754548270382321
$ws.Range("B13").Value = -0.2999999999999999
$ws.Range("C13").Value = 0.4
$ws.Range("D13").Value = -0.8
$ws.Range("E13").Value = 0.9
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = -0.20411440483222
$ws.Range("H13").Value = 0.2867771959663355
$ws.Range("B14").Value = -0.5
$ws.Range("C14").Value = 0.9
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = -0.2201833172501667
$ws.Range("F14").Value = 0.2499791286330684
$ws.Range("B15").Value = 0.09999999999999998
$ws.Range("C15").Value = -0.2000000000000001
$ws.Range("D15").Value = 0.1999999999999999
$ws.Range("E15").Value = 0.09999999999999991
$ws.Range("F15").Value = -0.7000000000000001
$ws.Range("G15").Value = 0.3619385045446099
$ws.Range("H15").Value = 0.4440258437269816
$ws.Range("B16").Value = 0.09999999999999998
$ws.Range("C16").Value = 0.3
$ws.Range("D16").Value = -0.6
$ws.Range("E16").Value = 0.37078400593656
$ws.Range("F16").Value = 0.4291293167667899
$ws.Range("B17").Value = -0.3000000000000001
$ws.Range("C17").Value = 0.4
$ws.Range("D17").Value = 0.6
$ws.Range("E17").Value = 0.2
$ws.Range("F17").Value = -0.5
$ws.Range("G17").Value = 0.352030371842539
$ws.Range("H17").Value = 0.2025233107720449
$ws.Range("B18").Value = 0.1959171462662728
$ws.Range("C18").Value = 0.1
$ws.Range("D18").Value = -0.4
$ws.Range("E18").Value = 0.4
$ws.Range("F18").Value = 0.2
$ws.Range("G18").Value = -0.4
$ws.Range("H18").Value = 0.4
$ws.Range("I18").Value = -0.382481144815657
$ws.Range("J18").Value = -0.479917721171199
$ws.Range("B19").Value = -0.4
$ws.Range("C19").Value = 0.4
$ws.Range("D19").Value = 0.2
$ws.Range("E19").Value = -0.4
$ws.Range("F19").Value = 0.4
$ws.Range("G19").Value = -0.4
$ws.Range("H19").Value = -0.5
$ws.Range("B20").Value = 0.184309031391911
$ws.Range("C20").Value = -0.3490319194399757
$ws.Range("D20").Value = 0.4441446536451586
$ws.Range("E20").Value = -0.3696970991840849
$ws.Range("F20").Value = -0.4881520669947828
$ws.Range("G20").Value = 0.1232798963004633
$ws.Range("H20").Value = 0.4234449690039844
$ws.Range("I20").Value = 1.032457285035994
$ws.Range("J20").Value = 0.01048914795523898
$ws.Range("B21").Value = 0.6151519277310885
$ws.Range("C21").Value = -0.2931763466948021
$ws.Range("D21").Value = -0.5003364362792013
$ws.Range("E21").Value = 0.1329202529759511
$ws.Range("F21").Value = 0.425754494287973
$ws.Range("G21").Value = 1.034884731917155
$ws.Range("H21").Value = 0.01372976374571738
$ws.Range("B22").Value = -0.4556810254602444
$ws.Range("C22").Value = 0.3138655158257279
$ws.Range("D22").Value = 0.5144278353286429
$ws.Range("E22").Value = 1.036505425522674
$ws.Range("F22").Value = 0.03243477352088853
$ws.Range("G22").Value = -0.1646508830073682
$ws.Range("H22").Value = 0.5401474412755328
$ws.Range("I22").Value = -0.4515687328697356
$ws.Range("J22").Value = 0.62835724284007
$ws.Range("B23").Value = 0.191234484381759
$ws.Range("C23").Value = 0.524389953244618
$ws.Range("D23").Value = 1.040034915326271
$ws.Range("E23").Value = 0.0297951016611645
$ws.Range("F23").Value = -0.1617391303018729
$ws.Range("G23").Value = 0.5425821084835898
$ws.Range("H23").Value = -0.4493783636021329
$ws.Range("I23").Value = 0.6310192074810944
$ws.Range("B24").Value = 0.1283399026883052
$ws.Range("C24").Value = 1.069946439435145
$ws.Range("D24").Value = 0.03016339292225023
$ws.Range("E24").Value = -0.1761847410291747
$ws.Range("F24").Value = 0.5259379121631296
$ws.Range("G24").Value = -0.4706783544374164
$ws.Range("H24").Value = 0.6096460097612491
$ws.Range("B25").Value = 1.228981402283536
$ws.Range("C25").Value = 0.1467255298080423
$ws.Range("D25").Value = -0.1774509833985318
$ws.Range("E25").Value = 0.5229500875209663
$ws.Range("F25").Value = -0.4843494915527642
$ws.Range("G25").Value = 0.6017739134081601
$ws.Range("H25").Value = 0.7965668651186684
$ws.Range("I25").Value = 2.798161589704257
$ws.Range("B26").Value = 0.2
$ws.Range("C26").Value = -0.1
$ws.Range("D26").Value = 0.5
$ws.Range("E26").Value = -0.5
$ws.Range("F26").Value = 0.6
$ws.Range("G26").Value = 0.8
$ws.Range("H26").Value = 2.8
$ws.Range("B27").Value = -0.3781284012960198
$ws.Range("C27").Value = 0.6143801122861683
$ws.Range("D27").Value = -0.4809023813400933
$ws.Range("E27").Value = 0.5922977141344893
$ws.Range("F27").Value = 0.7900134395145342
$ws.Range("G27").Value = 2.793803210748177
$ws.Range("B28").Value = 0.2
$ws.Range("C28").Value = -0.4
$ws.Range("D28").Value = 0.6
$ws.Range("E28").Value = 0.8
$ws.Range("F28").Value = 2.8
$ws.Range("B29").Value = -0.4349450679668658
$ws.Range("C29").Value = 0.6436351713107291
$ws.Range("D29").Value = 0.7850103128411116
$ws.Range("E29").Value = 2.793691431320905
$ws.Range("B30").Value = 0.7023597690241736
$ws.Range("C30").Value = 0.7403290071467026
$ws.Range("D30").Value = 2.796086438208883
$ws.Range("B31").Value = 0.2574694337905873
$ws.Range("C31").Value = 2.204570906693543
$ws.Range("B32").Value = -2.600000000000001
